# Add "Clinton Votes" (col F) and "Trump Votes" (col G) election-results
# columns to the state pay/strike table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("F1").Value = "Clinton Votes"
$ws.Range("G1").Value = "Trump Votes"

# Per-row vote totals, in the same row order (2..52) as the existing
# alphabetically-sorted state rows already on the sheet.
$clintonVotes = @(718084, 93007, 936250, 378729, 7362490, 1208095, 884432, 235581, 560223, 4485745, 1837300, 251853, 189677, 2977498, 1031953, 650790, 414788, 628834, 779535, 354873, 1497951, 1964768, 2268193, 1366676, 462001, 1054889, 175521, 273858, 537753, 348521, 2021756, 380724, 4143874, 2162074, 216133, 2317001, 419788, 934631, 2844705, 249902, 849469, 177442, 867110, 3867816, 274188, 178179, 1916845, 1610524, 187457, 1382210, 55949)
$trumpVotes   = @(1306925, 130415, 1021154, 677904, 3916209, 1136354, 668266, 185103, 11553, 4605515, 2068623, 121648, 407199, 2118179, 1556220, 798923, 656009, 1202942, 1178004, 334838, 873646, 1083069, 2279805, 1322891, 678457, 1585753, 274120, 485819, 511319, 345789, 1535513, 315875, 2640570, 2339603, 93526, 2771984, 947934, 742506, 2912941, 179421, 1143611, 227701, 1517402, 4681590, 452086, 95053, 1731156, 1129120, 486198, 1409467, 174248)

for ($i = 0; $i -lt $clintonVotes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $clintonVotes[$i]
    $ws.Cells.Item($row, 7).Value = $trumpVotes[$i]
}

# Match the saved view state: scrolled down with G53 (first empty cell
# below the new data) selected.
$ws.Range("G53").Select()
